$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get plain decimal-looking strings (e.g. "1.00", "118.60").
# Pre-format them as Text so Excel keeps the exact display string (incl. trailing
# zeros) instead of silently coercing to a number (which would drop them).
$textCells = @("D4", "D5", "D6", "D9", "D14", "D19", "D20", "D21", "D22", "D23", "D25", "D28", "D31", "D32", "D33", "D35", "D36", "D37", "D39", "D41", "D43", "D44", "D47", "D49")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# --- Cell value updates (per upstream cryptos-list refresh) ---
$ws.Range("D2").Value = "64.800.44"
$ws.Range("E2").Value = "  +5.20%  "
$ws.Range("D3").Value = "3.099.40"
$ws.Range("E3").Value = "  +3.11%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "558.43"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "143.56"
$ws.Range("E6").Value = "  +9.78%  "
$ws.Range("D8").Value = "3.095.82"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  +17.45%  "
$ws.Range("E11").Value = "  +4.87%  "
$ws.Range("E12").Value = "  +3.74%  "
$ws.Range("E13").Value = "  +4.45%  "
$ws.Range("D14").Value = "35.30"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").Value = "3.605.60"
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("D16").Value = "64.787.60"
$ws.Range("E16").Value = "  +5.06%  "
$ws.Range("D17").Value = "3.101.26"
$ws.Range("E17").Value = "  +3.36%  "
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "6.82"
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("D20").Value = "483.88"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "13.82"
$ws.Range("E21").Value = "  +4.76%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "0.675"
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "7.59"
$ws.Range("E23").Value = "  +8.72%  "
$ws.Range("E24").Value = "  +10.82%  "
$ws.Range("D25").Value = "80.88"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("D28").Value = "8.17"
$ws.Range("E28").Value = "  +6.85%  "
$ws.Range("E29").Value = "  +8.38%  "
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "26.06"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("D33").Value = "2.45"
$ws.Range("E33").Value = "  +5.04%  "
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("D35").Value = "6.20"
$ws.Range("E35").Value = "  +5.90%  "
$ws.Range("D36").Value = "54.94"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "466.05"
$ws.Range("E37").Value = "  +5.78%  "
$ws.Range("E38").Value = "  +6.66%  "
$ws.Range("D39").Value = "0.0827"
$ws.Range("E39").Value = "  +4.23%  "
$ws.Range("D40").Value = "3.015.01"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  +15.56%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").Value = "8.26"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").Value = "28.44"
$ws.Range("E44").Value = "  +9.51%  "
$ws.Range("E45").Value = "  +7.71%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  +8.70%  "
$ws.Range("E48").Value = "  +4.09%  "
$ws.Range("D49").Value = "118.60"
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("D50").Value = "0.0₃0515"
$ws.Range("E50").Value = "  +6.65%  "
$ws.Range("E51").Value = "  +2.84%  "
